$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The "date" and "time" types are no longer used for these fields;
# change them to "text" so the underlying shared strings "date"/"time"
# become unreferenced and are dropped on save.
$ws.Range("C2").Value = "text"
$ws.Range("C4").Value = "text"
$ws.Range("C6").Value = "text"
